$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place run edits on shared strings) ---
$ws.Range("A8").Characters(21, 2).Text = "43"
$ws.Range("C9").Characters(27, 10).Text = "10/24/2022"
$ws.Range("C9").Characters(48, 10).Text = "10/30/2022"

# Row 15
$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 90
$ws.Range("L15").Value = 5.555555555555
$ws.Range("M15").Value = -5
$ws.Range("N15").Value = -71.212121212121

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -66.666666666666
$ws.Range("I16").Value = 133
$ws.Range("J16").Value = 145
$ws.Range("K16").Value = -8.275862068965
$ws.Range("L16").Value = -13.071895424836
$ws.Range("M16").Value = -53.003533568904
$ws.Range("N16").Value = -86.922320550639

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -36.363636363636
$ws.Range("I17").Value = 281
$ws.Range("J17").Value = 280
$ws.Range("K17").Value = 0.357142857142
$ws.Range("L17").Value = 15.163934426229
$ws.Range("M17").Value = 3.308823529411
$ws.Range("N17").Value = -62.231182795698

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 22
$ws.Range("H18").Value = 37.5
$ws.Range("I18").Value = 170
$ws.Range("J18").Value = 141
$ws.Range("K18").Value = 20.567375886524
$ws.Range("L18").Value = -18.660287081339
$ws.Range("M18").Value = -13.705583756345
$ws.Range("N18").Value = -74.436090225563

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 5.714285714285
$ws.Range("I19").Value = 358
$ws.Range("J19").Value = 309
$ws.Range("K19").Value = 15.857605177993
$ws.Range("L19").Value = 10.835913312693
$ws.Range("M19").Value = 24.738675958188
$ws.Range("N19").Value = 16.612377850162

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 116
$ws.Range("J20").Value = 84
$ws.Range("K20").Value = 38.095238095238
$ws.Range("L20").Value = 27.472527472527
$ws.Range("M20").Value = 58.904109589041
$ws.Range("N20").Value = -77.075098814229

# Row 21
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 111
$ws.Range("H21").Value = -11.711711711711
$ws.Range("I21").Value = 1082
$ws.Range("J21").Value = 978
$ws.Range("K21").Value = 10.633946830265
$ws.Range("L21").Value = 3.441682600382
$ws.Range("M21").Value = -5.74912891986
$ws.Range("N21").Value = -67.48798076923

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("L22").Value = 8.333333333333

# Row 23
$ws.Range("I14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 72
$ws.Range("J23").Value = 92
$ws.Range("K23").Value = -21.739130434782
$ws.Range("L23").Value = -7.692307692307
$ws.Range("M23").Value = 7.462686567164

# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 157.142857142857
$ws.Range("F24").Value = 87
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = 117.5
$ws.Range("I24").Value = 746
$ws.Range("J24").Value = 491
$ws.Range("K24").Value = 51.93482688391
$ws.Range("L24").Value = 10.355029585798
$ws.Range("M24").Value = 9.705882352941

# Row 25
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -32
$ws.Range("I25").Value = 335
$ws.Range("J25").Value = 285
$ws.Range("K25").Value = 17.543859649122
$ws.Range("L25").Value = 9.83606557377
$ws.Range("M25").Value = -51.659451659451

# Row 26
$ws.Range("I14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 25
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = 8.695652173913
$ws.Range("L26").Value = 0

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 27
$ws.Range("K27").Value = -34.146341463414
$ws.Range("L27").Value = -41.304347826087

# Row 28
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -2.777777777777
$ws.Range("L28").Value = -27.083333333333
$ws.Range("N28").Value = -77.987421383647

# Row 29
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 27
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -30.76923076923
$ws.Range("N29").Value = -81.118881118881

# Row 30
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
